$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "contact": add an ORCID column (between phone and linkedin) and
# update the institute affiliation (EGID - UMR 8199 -> EGID - UMR 1283).
# ---------------------------------------------------------------------------
$contact = $wb.Worksheets.Item("contact")

# Insert a new column before the current "linkedin" column (F) for "orcid".
$contact.Columns("F").Insert()

$contact.Range("F1").Value = "orcid"
$contact.Range("F2").Value = "0000-0002-3396-4549"

# Update the institute / affiliation text for the contact block.
$contact.Range("B2").Value = "[EGID - UMR 1283](http://www.good.cnrs.fr/?lang=en)"

$contact.Columns("F").ColumnWidth = 18.5

# ---------------------------------------------------------------------------
# Sheet "experience": close out the previous position (end date) and add the
# new position at the renamed institute (EGID - UMR 1283).
# ---------------------------------------------------------------------------
$experience = $wb.Worksheets.Item("experience")

# The "Head of the Biostatistic Team" role at EGID - UMR 8199 now ends
# Dec. 2019 instead of being the ongoing ("Present") position.
$experience.Range("E6").Value = "Dec. 2019"

# Add the new row describing the continuation of the role at the renamed
# institute, starting Jan. 2020 and still ongoing ("Present").
$experience.Range("A7").Value = "Head of the Biostatistic Team"
$experience.Range("B7").Value = "Functional (Epi)genomics and Molecular  `nPhysiology of Diabetes and Associated Diseases  `nEGID - UMR 1283  `n(European Genomics Institute for Diabetes)"
$experience.Range("B7").WrapText = $true
$experience.Range("C7").Value = "Lille, France"
$experience.Range("D7").Value = "Jan. 2020"
$experience.Range("E7").Value = "Present"
$experience.Range("F7").Value = "Genome-wide association studies, experimental design, -omics data analysis, methodological developments, team management"

$experience.Rows("7").RowHeight = 75

$experience.Range("B7").Select()

# Keep "contact" as the active/selected sheet and tab when the workbook is
# saved (matching the original workbook's active tab).
$contact.Activate()
$contact.Range("A1").Select()
